# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows for "Provincia de Melipilla" /
# Frutilla just above the existing row 273 block. Every row below shifts
# down by two (old row N -> new row N+2); the sheet's used range grows
# from A1:T338 to A1:T340.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 273. Excel pushes
# everything from the old row 273 onward down by two rows, preserving
# the per-row formatting (date style on column D, etc.) via InsertCopy
# semantics.
$ws.Rows.Item(273).Insert()
$ws.Rows.Item(273).Insert()

# --- New row 273: "Especial" tier record replaced by "Primera" tier for
#     Provincia de Melipilla, week of 2022-10-?? (serial 44841) ---
$ws.Cells.Item(273,1).Value2  = 7
$ws.Cells.Item(273,2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(273,3).Value2  = "Ñuble"
$ws.Cells.Item(273,4).Value2  = 44841
$ws.Cells.Item(273,5).Value2  = 16
$ws.Cells.Item(273,6).Value2  = "Fruta"
$ws.Cells.Item(273,7).Value2  = 100101
$ws.Cells.Item(273,8).Value2  = "Berries"
$ws.Cells.Item(273,9).Value2  = 100112025
$ws.Cells.Item(273,10).Value2 = "Frutilla"
$ws.Cells.Item(273,11).Value2 = "Sin especificar"
$ws.Cells.Item(273,12).Value2 = "Primera"
$ws.Cells.Item(273,13).Value2 = 120
$ws.Cells.Item(273,14).Value2 = 13000
$ws.Cells.Item(273,15).Value2 = 14000
$ws.Cells.Item(273,16).Value2 = 13500
$ws.Cells.Item(273,17).Value2 = "$/bandeja 7 kilos"
$ws.Cells.Item(273,18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(273,19).Value2 = 1929
$ws.Cells.Item(273,20).Value2 = 7

# --- New row 274: "Segunda" tier record for Provincia de Melipilla,
#     same week (serial 44841) ---
$ws.Cells.Item(274,1).Value2  = 7
$ws.Cells.Item(274,2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(274,3).Value2  = "Ñuble"
$ws.Cells.Item(274,4).Value2  = 44841
$ws.Cells.Item(274,5).Value2  = 16
$ws.Cells.Item(274,6).Value2  = "Fruta"
$ws.Cells.Item(274,7).Value2  = 100101
$ws.Cells.Item(274,8).Value2  = "Berries"
$ws.Cells.Item(274,9).Value2  = 100112025
$ws.Cells.Item(274,10).Value2 = "Frutilla"
$ws.Cells.Item(274,11).Value2 = "Sin especificar"
$ws.Cells.Item(274,12).Value2 = "Segunda"
$ws.Cells.Item(274,13).Value2 = 60
$ws.Cells.Item(274,14).Value2 = 9000
$ws.Cells.Item(274,15).Value2 = 9000
$ws.Cells.Item(274,16).Value2 = 9000
$ws.Cells.Item(274,17).Value2 = "$/bandeja 7 kilos"
$ws.Cells.Item(274,18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(274,19).Value2 = 1286
$ws.Cells.Item(274,20).Value2 = 7
